$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") and column E ("Volume(1h)") hold plain text in this
# sheet (e.g. "62.871.69", "  +1.33%  "). A leading apostrophe forces
# Excel to keep a numeric-looking string (e.g. "571.30") as literal text
# instead of silently converting it to the number 571.3.
$ws.Range("D2").Value = "62.871.69"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "2.437.78"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'571.30"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").Value = "'146.58"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").Value = "'0.112"
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("D12").Value = "'0.358"
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").Value = "'26.97"
$ws.Range("E13").Value = "  +5.09%  "
$ws.Range("D14").Value = "'0.0000182"
$ws.Range("E14").Value = "  +4.35%  "
$ws.Range("D15").Value = "2.876.60"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "62.601.84"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "2.441.36"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "'11.30"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").Value = "'7.08"
$ws.Range("E19").Value = "  +3.84%  "
$ws.Range("D20").Value = "'325.30"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'1.85"
$ws.Range("E23").Value = "  +5.27%  "
$ws.Range("D24").Value = "'67.28"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").Value = "'628.34"
$ws.Range("E25").Value = "  +10.59%  "
$ws.Range("D26").Value = "'8.70"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").Value = "'0.0000103"
$ws.Range("E27").Value = "  +10.12%  "
$ws.Range("D28").Value = "2.556.44"
$ws.Range("D29").Value = "'8.51"
$ws.Range("E29").Value = "  +3.77%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E31").Value = "  +4.83%  "
$ws.Range("E32").Value = "  -3.34%  "
$ws.Range("D33").Value = "'1.89"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "'1.52"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").Value = "'4.95"
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("D38").Value = "'5.45"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "'18.80"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("D40").Value = "'1.84"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").Value = "'148.28"
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("D42").Value = "'2.58"
$ws.Range("E42").Value = "  +14.90%  "
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "'150.38"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").Value = "'3.70"
$ws.Range("E45").Value = "  +2.49%  "
$ws.Range("D46").Value = "'0.0541"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").Value = "'20.84"
$ws.Range("E47").Value = "  +4.74%  "
$ws.Range("D48").Value = "'0.606"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("D50").Value = "'0.0925"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("E51").Value = "  +4.47%  "
